$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2343, 1).Value = 4
$ws.Cells.Item(2343, 2).Value = 3
$ws.Cells.Item(2343, 3).Value = 3
$ws.Cells.Item(2343, 4).Value = 0
$ws.Cells.Item(2344, 1).Value = 5
$ws.Cells.Item(2344, 2).Value = 2
$ws.Cells.Item(2344, 3).Value = 5
$ws.Cells.Item(2344, 4).Value = 0
$ws.Cells.Item(2345, 1).Value = 7
$ws.Cells.Item(2345, 2).Value = 1
$ws.Cells.Item(2345, 3).Value = 6
$ws.Cells.Item(2345, 4).Value = 2
$ws.Cells.Item(2346, 1).Value = 5
$ws.Cells.Item(2346, 2).Value = 2
$ws.Cells.Item(2346, 3).Value = 5
$ws.Cells.Item(2346, 4).Value = 1
$ws.Cells.Item(2347, 1).Value = 2
$ws.Cells.Item(2347, 2).Value = 2
$ws.Cells.Item(2347, 3).Value = 3
$ws.Cells.Item(2347, 4).Value = 1
$ws.Cells.Item(2348, 1).Value = 5
$ws.Cells.Item(2348, 2).Value = 1
$ws.Cells.Item(2348, 3).Value = 5
$ws.Cells.Item(2348, 4).Value = 2
$ws.Cells.Item(2349, 1).Value = 5
$ws.Cells.Item(2349, 2).Value = 2
$ws.Cells.Item(2349, 3).Value = 5
$ws.Cells.Item(2349, 4).Value = 0
$ws.Cells.Item(2350, 1).Value = 4
$ws.Cells.Item(2350, 2).Value = 0
$ws.Cells.Item(2350, 3).Value = 4
$ws.Cells.Item(2350, 4).Value = 2
$ws.Cells.Item(2351, 1).Value = 3
$ws.Cells.Item(2351, 2).Value = 2
$ws.Cells.Item(2351, 3).Value = 4
$ws.Cells.Item(2351, 4).Value = 0
$ws.Cells.Item(2352, 1).Value = 6
$ws.Cells.Item(2352, 2).Value = 0
$ws.Cells.Item(2352, 3).Value = 5
$ws.Cells.Item(2352, 4).Value = 3
$ws.Cells.Item(2353, 1).Value = 5
$ws.Cells.Item(2353, 2).Value = 2
$ws.Cells.Item(2353, 3).Value = 6
$ws.Cells.Item(2353, 4).Value = 1
$ws.Cells.Item(2354, 1).Value = 3
$ws.Cells.Item(2354, 2).Value = 0
$ws.Cells.Item(2354, 3).Value = 3
$ws.Cells.Item(2354, 4).Value = 3
$ws.Cells.Item(2355, 1).Value = 2
$ws.Cells.Item(2355, 2).Value = 2
$ws.Cells.Item(2355, 3).Value = 3
$ws.Cells.Item(2355, 4).Value = 1
$ws.Cells.Item(2356, 1).Value = 7
$ws.Cells.Item(2356, 2).Value = 2
$ws.Cells.Item(2356, 3).Value = 6
$ws.Cells.Item(2356, 4).Value = 0
$ws.Cells.Item(2357, 1).Value = 4
$ws.Cells.Item(2357, 2).Value = 1
$ws.Cells.Item(2357, 3).Value = 5
$ws.Cells.Item(2357, 4).Value = 2
$ws.Cells.Item(2358, 1).Value = 4
$ws.Cells.Item(2358, 2).Value = 2
$ws.Cells.Item(2358, 3).Value = 4
$ws.Cells.Item(2358, 4).Value = 0
$ws.Cells.Item(2359, 1).Value = 4
$ws.Cells.Item(2359, 2).Value = 2
$ws.Cells.Item(2359, 3).Value = 4
$ws.Cells.Item(2359, 4).Value = 0
$ws.Cells.Item(2360, 1).Value = 6
$ws.Cells.Item(2360, 2).Value = 2
$ws.Cells.Item(2360, 3).Value = 6
$ws.Cells.Item(2360, 4).Value = 1
$ws.Cells.Item(2361, 1).Value = 3
$ws.Cells.Item(2361, 2).Value = 1
$ws.Cells.Item(2361, 3).Value = 4
$ws.Cells.Item(2361, 4).Value = 2
$ws.Cells.Item(2362, 1).Value = 6
$ws.Cells.Item(2362, 2).Value = 2
$ws.Cells.Item(2362, 3).Value = 5
$ws.Cells.Item(2362, 4).Value = 0
$ws.Cells.Item(2363, 1).Value = 3
$ws.Cells.Item(2363, 2).Value = 0
$ws.Cells.Item(2363, 3).Value = 3
$ws.Cells.Item(2363, 4).Value = 3
$ws.Cells.Item(2364, 1).Value = 5
$ws.Cells.Item(2364, 2).Value = 2
$ws.Cells.Item(2364, 3).Value = 4
$ws.Cells.Item(2364, 4).Value = 1
$ws.Cells.Item(2365, 1).Value = 4
$ws.Cells.Item(2365, 2).Value = 2
$ws.Cells.Item(2365, 3).Value = 5
$ws.Cells.Item(2365, 4).Value = 0
$ws.Cells.Item(2366, 1).Value = 5
$ws.Cells.Item(2366, 2).Value = 2
$ws.Cells.Item(2366, 3).Value = 6
$ws.Cells.Item(2366, 4).Value = 0
$ws.Cells.Item(2367, 1).Value = 4
$ws.Cells.Item(2367, 2).Value = 0
$ws.Cells.Item(2367, 3).Value = 5
$ws.Cells.Item(2367, 4).Value = 2
$ws.Cells.Item(2368, 1).Value = 5
$ws.Cells.Item(2368, 2).Value = 1
$ws.Cells.Item(2368, 3).Value = 7
$ws.Cells.Item(2368, 4).Value = 2
$ws.Cells.Item(2369, 1).Value = 3
$ws.Cells.Item(2369, 2).Value = 0
$ws.Cells.Item(2369, 3).Value = 3
$ws.Cells.Item(2369, 4).Value = 3
$ws.Cells.Item(2370, 1).Value = 4
$ws.Cells.Item(2370, 2).Value = 0
$ws.Cells.Item(2370, 3).Value = 3
$ws.Cells.Item(2370, 4).Value = 2
$ws.Cells.Item(2371, 1).Value = 6
$ws.Cells.Item(2371, 2).Value = 3
$ws.Cells.Item(2371, 3).Value = 6
$ws.Cells.Item(2371, 4).Value = 0
$ws.Cells.Item(2372, 1).Value = 4
$ws.Cells.Item(2372, 2).Value = 0
$ws.Cells.Item(2372, 3).Value = 3
$ws.Cells.Item(2372, 4).Value = 2
$ws.Cells.Item(2373, 1).Value = 4
$ws.Cells.Item(2373, 2).Value = 2
$ws.Cells.Item(2373, 3).Value = 4
$ws.Cells.Item(2373, 4).Value = 1
$ws.Cells.Item(2374, 1).Value = 3
$ws.Cells.Item(2374, 2).Value = 3
$ws.Cells.Item(2374, 3).Value = 2
$ws.Cells.Item(2374, 4).Value = 0
$ws.Cells.Item(2375, 1).Value = 6
$ws.Cells.Item(2375, 2).Value = 0
$ws.Cells.Item(2375, 3).Value = 5
$ws.Cells.Item(2375, 4).Value = 2
$ws.Cells.Item(2376, 1).Value = 3
$ws.Cells.Item(2376, 2).Value = 0
$ws.Cells.Item(2376, 3).Value = 3
$ws.Cells.Item(2376, 4).Value = 3
$ws.Cells.Item(2377, 1).Value = 3
$ws.Cells.Item(2377, 2).Value = 2
$ws.Cells.Item(2377, 3).Value = 6
$ws.Cells.Item(2377, 4).Value = 1
$ws.Cells.Item(2378, 1).Value = 4
$ws.Cells.Item(2378, 2).Value = 1
$ws.Cells.Item(2378, 3).Value = 5
$ws.Cells.Item(2378, 4).Value = 2
$ws.Cells.Item(2379, 1).Value = 4
$ws.Cells.Item(2379, 2).Value = 2
$ws.Cells.Item(2379, 3).Value = 5
$ws.Cells.Item(2379, 4).Value = 0
$ws.Cells.Item(2380, 1).Value = 3
$ws.Cells.Item(2380, 2).Value = 3
$ws.Cells.Item(2380, 3).Value = 4
$ws.Cells.Item(2380, 4).Value = 0
$ws.Cells.Item(2381, 1).Value = 5
$ws.Cells.Item(2381, 2).Value = 2
$ws.Cells.Item(2381, 3).Value = 3
$ws.Cells.Item(2381, 4).Value = 1
$ws.Cells.Item(2382, 1).Value = 4
$ws.Cells.Item(2382, 2).Value = 2
$ws.Cells.Item(2382, 3).Value = 3
$ws.Cells.Item(2382, 4).Value = 1
$ws.Cells.Item(2383, 1).Value = 4
$ws.Cells.Item(2383, 2).Value = 1
$ws.Cells.Item(2383, 3).Value = 5
$ws.Cells.Item(2383, 4).Value = 2
$ws.Cells.Item(2384, 1).Value = 4
$ws.Cells.Item(2384, 2).Value = 1
$ws.Cells.Item(2384, 3).Value = 5
$ws.Cells.Item(2384, 4).Value = 2
$ws.Cells.Item(2385, 1).Value = 4
$ws.Cells.Item(2385, 2).Value = 3
$ws.Cells.Item(2385, 3).Value = 3
$ws.Cells.Item(2385, 4).Value = 0
$ws.Cells.Item(2386, 1).Value = 5
$ws.Cells.Item(2386, 2).Value = 0
$ws.Cells.Item(2386, 3).Value = 7
$ws.Cells.Item(2386, 4).Value = 2
$ws.Cells.Item(2387, 1).Value = 4
$ws.Cells.Item(2387, 2).Value = 3
$ws.Cells.Item(2387, 3).Value = 3
$ws.Cells.Item(2387, 4).Value = 0
$ws.Cells.Item(2388, 1).Value = 5
$ws.Cells.Item(2388, 2).Value = 1
$ws.Cells.Item(2388, 3).Value = 5
$ws.Cells.Item(2388, 4).Value = 2
$ws.Cells.Item(2389, 1).Value = 4
$ws.Cells.Item(2389, 2).Value = 1
$ws.Cells.Item(2389, 3).Value = 4
$ws.Cells.Item(2389, 4).Value = 2
$ws.Cells.Item(2390, 1).Value = 5
$ws.Cells.Item(2390, 2).Value = 0
$ws.Cells.Item(2390, 3).Value = 7
$ws.Cells.Item(2390, 4).Value = 3
$ws.Cells.Item(2391, 1).Value = 5
$ws.Cells.Item(2391, 2).Value = 2
$ws.Cells.Item(2391, 3).Value = 4
$ws.Cells.Item(2391, 4).Value = 1
$ws.Cells.Item(2392, 1).Value = 3
$ws.Cells.Item(2392, 2).Value = 1
$ws.Cells.Item(2392, 3).Value = 3
$ws.Cells.Item(2392, 4).Value = 2
$ws.Cells.Item(2393, 1).Value = 3
$ws.Cells.Item(2393, 2).Value = 0
$ws.Cells.Item(2393, 3).Value = 3
$ws.Cells.Item(2393, 4).Value = 3
$ws.Cells.Item(2394, 1).Value = 4
$ws.Cells.Item(2394, 2).Value = 2
$ws.Cells.Item(2394, 3).Value = 3
$ws.Cells.Item(2394, 4).Value = 1
$ws.Cells.Item(2395, 1).Value = 6
$ws.Cells.Item(2395, 2).Value = 1
$ws.Cells.Item(2395, 3).Value = 6
$ws.Cells.Item(2395, 4).Value = 2
$ws.Cells.Item(2396, 1).Value = 5
$ws.Cells.Item(2396, 2).Value = 2
$ws.Cells.Item(2396, 3).Value = 5
$ws.Cells.Item(2396, 4).Value = 0
$ws.Cells.Item(2397, 1).Value = 3
$ws.Cells.Item(2397, 2).Value = 2
$ws.Cells.Item(2397, 3).Value = 2
$ws.Cells.Item(2397, 4).Value = 1
$ws.Cells.Item(2398, 1).Value = 5
$ws.Cells.Item(2398, 2).Value = 0
$ws.Cells.Item(2398, 3).Value = 5
$ws.Cells.Item(2398, 4).Value = 2
$ws.Cells.Item(2399, 1).Value = 5
$ws.Cells.Item(2399, 2).Value = 2
$ws.Cells.Item(2399, 3).Value = 5
$ws.Cells.Item(2399, 4).Value = 1
$ws.Cells.Item(2400, 1).Value = 4
$ws.Cells.Item(2400, 2).Value = 1
$ws.Cells.Item(2400, 3).Value = 3
$ws.Cells.Item(2400, 4).Value = 2
$ws.Cells.Item(2401, 1).Value = 6
$ws.Cells.Item(2401, 2).Value = 2
$ws.Cells.Item(2401, 3).Value = 5
$ws.Cells.Item(2401, 4).Value = 0
$ws.Cells.Item(2402, 1).Value = 5
$ws.Cells.Item(2402, 2).Value = 2
$ws.Cells.Item(2402, 3).Value = 4
$ws.Cells.Item(2402, 4).Value = 1
$ws.Cells.Item(2403, 1).Value = 3
$ws.Cells.Item(2403, 2).Value = 0
$ws.Cells.Item(2403, 3).Value = 4
$ws.Cells.Item(2403, 4).Value = 3
$ws.Cells.Item(2404, 1).Value = 4
$ws.Cells.Item(2404, 2).Value = 1
$ws.Cells.Item(2404, 3).Value = 4
$ws.Cells.Item(2404, 4).Value = 2
$ws.Cells.Item(2405, 1).Value = 5
$ws.Cells.Item(2405, 2).Value = 2
$ws.Cells.Item(2405, 3).Value = 6
$ws.Cells.Item(2405, 4).Value = 0
$ws.Cells.Item(2406, 1).Value = 4
$ws.Cells.Item(2406, 2).Value = 0
$ws.Cells.Item(2406, 3).Value = 4
$ws.Cells.Item(2406, 4).Value = 2
$ws.Cells.Item(2407, 1).Value = 6
$ws.Cells.Item(2407, 2).Value = 0
$ws.Cells.Item(2407, 3).Value = 6
$ws.Cells.Item(2407, 4).Value = 3
$ws.Cells.Item(2408, 1).Value = 4
$ws.Cells.Item(2408, 2).Value = 0
$ws.Cells.Item(2408, 3).Value = 3
$ws.Cells.Item(2408, 4).Value = 2
$ws.Cells.Item(2409, 1).Value = 6
$ws.Cells.Item(2409, 2).Value = 2
$ws.Cells.Item(2409, 3).Value = 6
$ws.Cells.Item(2409, 4).Value = 0
$ws.Cells.Item(2410, 1).Value = 3
$ws.Cells.Item(2410, 2).Value = 2
$ws.Cells.Item(2410, 3).Value = 3
$ws.Cells.Item(2410, 4).Value = 1
$ws.Cells.Item(2411, 1).Value = 4
$ws.Cells.Item(2411, 2).Value = 2
$ws.Cells.Item(2411, 3).Value = 3
$ws.Cells.Item(2411, 4).Value = 1
$ws.Cells.Item(2412, 1).Value = 5
$ws.Cells.Item(2412, 2).Value = 0
$ws.Cells.Item(2412, 3).Value = 5
$ws.Cells.Item(2412, 4).Value = 2
$ws.Cells.Item(2413, 1).Value = 3
$ws.Cells.Item(2413, 2).Value = 2
$ws.Cells.Item(2413, 3).Value = 2
$ws.Cells.Item(2413, 4).Value = 1
$ws.Cells.Item(2414, 1).Value = 4
$ws.Cells.Item(2414, 2).Value = 2
$ws.Cells.Item(2414, 3).Value = 3
$ws.Cells.Item(2414, 4).Value = 1
$ws.Cells.Item(2415, 1).Value = 5
$ws.Cells.Item(2415, 2).Value = 2
$ws.Cells.Item(2415, 3).Value = 7
$ws.Cells.Item(2415, 4).Value = 1
$ws.Cells.Item(2416, 1).Value = 4
$ws.Cells.Item(2416, 2).Value = 0
$ws.Cells.Item(2416, 3).Value = 3
$ws.Cells.Item(2416, 4).Value = 3

[void]$ws.Range("A2417").Select()
